# Add season-record columns (Wins / Losses / Ties) to the roster table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (AC1, style index 1:
# bold font, thin border, centered/top alignment) onto the three new header
# cells so they match the rest of the header row exactly.
$headerSrc = $ws.Range("AC1")
$headerSrc.Copy()
$headerDst = $ws.Range("AD1:AF1")
$headerDst.PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(1, 30).Value = "Wins"
$ws.Cells.Item(1, 31).Value = "Losses"
$ws.Cells.Item(1, 32).Value = "Ties"

# Fill the season record for every player row with the team's record.
for ($r = 2; $r -le 46; $r++) {
    $ws.Cells.Item($r, 30).Value = 66
    $ws.Cells.Item($r, 31).Value = 96
    $ws.Cells.Item($r, 32).Value = 0
}
